$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.1141855918187531
$ws.Cells.Item(2, 4).Value = 0.3333169489289673
$ws.Cells.Item(2, 5).Value = 0.07228057393488996
$ws.Cells.Item(2, 6).Value = 8.81180779629409
$ws.Cells.Item(2, 7).Value = 0.00268980014758706
$ws.Cells.Item(2, 9).Value = 0.4039108783602199
$ws.Cells.Item(2, 10).Value = 0.06544983559139794
$ws.Cells.Item(2, 13).Value = 3.712919696322587
$ws.Cells.Item(2, 14).Value = 1.277554709393584
$ws.Cells.Item(3, 2).Value = 0.1008095624010394
$ws.Cells.Item(3, 4).Value = 0.2980405007020011
$ws.Cells.Item(3, 5).Value = 0.06280600970178796
$ws.Cells.Item(3, 6).Value = 8.576231426350375
$ws.Cells.Item(3, 7).Value = 0.002707523031837855
$ws.Cells.Item(3, 9).Value = 0.411782040084729
$ws.Cells.Item(3, 10).Value = 0.06456493952439146
$ws.Cells.Item(3, 13).Value = 3.495903041952886
$ws.Cells.Item(3, 14).Value = 1.285322665297556
$ws.Cells.Item(4, 2).Value = 0.09258336114528731
$ws.Cells.Item(4, 4).Value = 0.276718284951869
$ws.Cells.Item(4, 5).Value = 0.05700749745576417
$ws.Cells.Item(4, 6).Value = 8.440308497472017
$ws.Cells.Item(4, 7).Value = 0.002718903784300632
$ws.Cells.Item(4, 9).Value = 0.4168841013767013
$ws.Cells.Item(4, 10).Value = 0.06402602727505879
$ws.Cells.Item(4, 13).Value = 3.36546262853517
$ws.Cells.Item(4, 14).Value = 1.290708908355413
$ws.Cells.Item(5, 2).Value = 0.08922811169776423
$ws.Cells.Item(5, 4).Value = 0.2681086155417347
$ws.Cells.Item(5, 5).Value = 0.05464823394373752
$ws.Cells.Item(5, 6).Value = 8.387054535638157
$ws.Cells.Item(5, 7).Value = 0.002723667960942203
$ws.Cells.Item(5, 9).Value = 0.4190309441090818
$ws.Cells.Item(5, 10).Value = 0.06380748318707319
$ws.Cells.Item(5, 13).Value = 3.312988843516109
$ws.Cells.Item(5, 14).Value = 1.293059640577219
$ws.Cells.Item(6, 2).Value = 0.08867080332441901
$ws.Cells.Item(6, 4).Value = 0.2666836163687094
$ws.Cells.Item(6, 5).Value = 0.05425667137317447
$ws.Cells.Item(6, 6).Value = 8.378339109065649
$ws.Cells.Item(6, 7).Value = 0.002724466713400778
$ws.Cells.Item(6, 9).Value = 0.4193915163515944
$ws.Cells.Item(6, 10).Value = 0.06377125731506439
$ws.Cells.Item(6, 13).Value = 3.304316088841716
$ws.Cells.Item(6, 14).Value = 1.293459408924193
$ws.Cells.Item(7, 2).Value = 0.09253812280854845
$ws.Cells.Item(7, 4).Value = 0.2766018580915102
$ws.Cells.Item(7, 5).Value = 0.05697566611002003
$ws.Cells.Item(7, 6).Value = 8.43958172277641
$ws.Cells.Item(7, 7).Value = 0.002718967522448847
$ws.Cells.Item(7, 9).Value = 0.4169127802383796
$ws.Cells.Item(7, 10).Value = 0.06402307564910714
$ws.Cells.Item(7, 13).Value = 3.364752220060382
$ws.Cells.Item(7, 14).Value = 1.290739979373882
$ws.Cells.Item(8, 2).Value = 0.109576482731697
$ws.Cells.Item(8, 4).Value = 0.3210801045543974
$ws.Cells.Item(8, 5).Value = 0.06900915686149034
$ws.Cells.Item(8, 6).Value = 8.728735505745419
$ws.Cells.Item(8, 7).Value = 0.002695808070731884
$ws.Cells.Item(8, 9).Value = 0.4065690172184357
$ws.Cells.Item(8, 10).Value = 0.06514377927114978
$ws.Cells.Item(8, 13).Value = 3.637494276173868
$ws.Cells.Item(8, 14).Value = 1.280105536868675
$ws.Cells.Item(9, 2).Value = 0.1428706530693233
$ws.Cells.Item(9, 4).Value = 0.4112470926662297
$ws.Cells.Item(9, 5).Value = 0.09280826514971352
$ws.Cells.Item(9, 6).Value = 9.367700595510769
$ws.Cells.Item(9, 7).Value = 0.002654304278975268
$ws.Cells.Item(9, 9).Value = 0.3884186046752465
$ws.Cells.Item(9, 10).Value = 0.06737872431420655
$ws.Cells.Item(9, 13).Value = 4.195809051292827
$ws.Cells.Item(9, 14).Value = 1.264113668017842
$ws.Cells.Item(10, 2).Value = 0.1672453775322538
$ws.Cells.Item(10, 4).Value = 0.4796778710994545
$ws.Cells.Item(10, 5).Value = 0.110491996975739
$ws.Cells.Item(10, 6).Value = 9.885003994540796
$ws.Cells.Item(10, 7).Value = 0.002626130082598916
$ws.Cells.Item(10, 9).Value = 0.3763814195515653
$ws.Cells.Item(10, 10).Value = 0.06904683420804503
$ws.Cells.Item(10, 13).Value = 4.622067154923371
$ws.Cells.Item(10, 14).Value = 1.255288292216051
$ws.Cells.Item(11, 2).Value = 0.1783124145359238
$ws.Cells.Item(11, 4).Value = 0.5113734602943509
$ws.Cells.Item(11, 5).Value = 0.1185974515411843
$ws.Cells.Item(11, 6).Value = 10.13165864121345
$ws.Cells.Item(11, 7).Value = 0.002613801561328294
$ws.Cells.Item(11, 9).Value = 0.3711866250540261
$ws.Cells.Item(11, 10).Value = 0.06981217341703427
$ws.Cells.Item(11, 13).Value = 4.819881476264044
$ws.Cells.Item(11, 14).Value = 1.251899959915477
$ws.Cells.Item(12, 2).Value = 0.1824998850534314
$ws.Cells.Item(12, 4).Value = 0.5234649031838217
$ws.Cells.Item(12, 5).Value = 0.121677054413631
$ws.Cells.Item(12, 6).Value = 10.22676976099706
$ws.Cells.Item(12, 7).Value = 0.00260920202682191
$ws.Cells.Item(12, 9).Value = 0.3692598638240856
$ws.Cells.Item(12, 10).Value = 0.07010299105635909
$ws.Cells.Item(12, 13).Value = 4.895385557163195
$ws.Cells.Item(12, 14).Value = 1.250706275362205
$ws.Cells.Item(13, 2).Value = 0.181598193804632
$ws.Cells.Item(13, 4).Value = 0.5208567104844519
$ws.Cells.Item(13, 5).Value = 0.1210133270965343
$ws.Cells.Item(13, 6).Value = 10.2062085977023
$ws.Cells.Item(13, 7).Value = 0.002610189568728901
$ws.Cells.Item(13, 9).Value = 0.3696730299177027
$ws.Cells.Item(13, 10).Value = 0.07004031271225841
$ws.Cells.Item(13, 13).Value = 4.879097325965574
$ws.Cells.Item(13, 14).Value = 1.250959391027067
$ws.Cells.Item(14, 2).Value = 0.1786569898254129
$ws.Cells.Item(14, 4).Value = 0.512366401690258
$ws.Cells.Item(14, 5).Value = 0.1188505983628332
$ws.Cells.Item(14, 6).Value = 10.13944878499154
$ws.Cells.Item(14, 7).Value = 0.002613421778481135
$ws.Cells.Item(14, 9).Value = 0.3710272997330826
$ws.Cells.Item(14, 10).Value = 0.06983607872678732
$ws.Cells.Item(14, 13).Value = 4.826081093152368
$ws.Cells.Item(14, 14).Value = 1.251799966797847
$ws.Cells.Item(15, 2).Value = 0.1768549679028126
$ws.Cells.Item(15, 4).Value = 0.5071776667383574
$ws.Cells.Item(15, 5).Value = 0.1175272433493433
$ws.Cells.Item(15, 6).Value = 10.09878136584575
$ws.Cells.Item(15, 7).Value = 0.002615410552243876
$ws.Cells.Item(15, 9).Value = 0.3718620893955382
$ws.Cells.Item(15, 10).Value = 0.06971111176991229
$ws.Cells.Item(15, 13).Value = 4.793685796131001
$ws.Cells.Item(15, 14).Value = 1.252326466959715
$ws.Cells.Item(16, 2).Value = 0.1665216663427742
$ws.Cells.Item(16, 4).Value = 0.4776185225325662
$ws.Cells.Item(16, 5).Value = 0.1099636382679137
$ws.Cells.Item(16, 6).Value = 9.869119198261899
$ws.Cells.Item(16, 7).Value = 0.002626945501876327
$ws.Cells.Item(16, 9).Value = 0.3767265647750762
$ws.Cells.Item(16, 10).Value = 0.06899695385605042
$ws.Cells.Item(16, 13).Value = 4.60922103222461
$ws.Cells.Item(16, 14).Value = 1.255522274653885
$ws.Cells.Item(17, 2).Value = 0.1601768642966164
$ws.Cells.Item(17, 4).Value = 0.4596352260936101
$ws.Cells.Item(17, 5).Value = 0.1053402638266334
$ws.Cells.Item(17, 6).Value = 9.73118591536155
$ws.Cells.Item(17, 7).Value = 0.002634146001049131
$ws.Cells.Item(17, 9).Value = 0.3797827215258618
$ws.Cells.Item(17, 10).Value = 0.06856055240892545
$ws.Cells.Item(17, 13).Value = 4.497082403350447
$ws.Cells.Item(17, 14).Value = 1.257642750860427
$ws.Cells.Item(18, 2).Value = 0.1565255312482918
$ws.Cells.Item(18, 4).Value = 0.4493445143411918
$ws.Cells.Item(18, 5).Value = 0.1026866739667796
$ws.Cells.Item(18, 6).Value = 9.652911875354789
$ws.Cells.Item(18, 7).Value = 0.002638333549291659
$ws.Cells.Item(18, 9).Value = 0.3815669899135106
$ws.Cells.Item(18, 10).Value = 0.06831015425656517
$ws.Cells.Item(18, 13).Value = 4.432948852783539
$ws.Cells.Item(18, 14).Value = 1.258921389601824
$ws.Cells.Item(19, 2).Value = 0.1552889244266993
$ws.Cells.Item(19, 4).Value = 0.4458690998288262
$ws.Cells.Item(19, 5).Value = 0.1017891414722598
$ws.Cells.Item(19, 6).Value = 9.626589608630297
$ws.Cells.Item(19, 7).Value = 0.002639759318483953
$ws.Cells.Item(19, 9).Value = 0.3821756552940574
$ws.Cells.Item(19, 10).Value = 0.06822547634568821
$ws.Cells.Item(19, 13).Value = 4.411296118482426
$ws.Cells.Item(19, 14).Value = 1.259364467974805
$ws.Cells.Item(20, 2).Value = 0.1608524857101941
$ws.Cells.Item(20, 4).Value = 0.4615440597217457
$ws.Cells.Item(20, 5).Value = 0.1058318345024958
$ws.Cells.Item(20, 6).Value = 9.745758724454447
$ws.Cells.Item(20, 7).Value = 0.002633374742171574
$ws.Cells.Item(20, 9).Value = 0.3794546512457586
$ws.Cells.Item(20, 10).Value = 0.06860694466565675
$ws.Cells.Item(20, 13).Value = 4.508981669383701
$ws.Cells.Item(20, 14).Value = 1.257410921459552
$ws.Cells.Item(21, 2).Value = 0.1795209870099086
$ws.Cells.Item(21, 4).Value = 0.514857732877033
$ws.Cells.Item(21, 5).Value = 0.1194855541540534
$ws.Cells.Item(21, 6).Value = 10.15901074577329
$ws.Cells.Item(21, 7).Value = 0.00261247053690165
$ws.Cells.Item(21, 9).Value = 0.3706284216366491
$ws.Cells.Item(21, 10).Value = 0.06989603951516443
$ws.Cells.Item(21, 13).Value = 4.84163678602394
$ws.Cells.Item(21, 14).Value = 1.251550648232111
$ws.Cells.Item(22, 2).Value = 0.1917021000906232
$ws.Cells.Item(22, 4).Value = 0.5502239404304987
$ws.Cells.Item(22, 5).Value = 0.1284694905734867
$ws.Cells.Item(22, 6).Value = 10.43909211422579
$ws.Cells.Item(22, 7).Value = 0.002599210094606625
$ws.Cells.Item(22, 9).Value = 0.3650954174035519
$ws.Cells.Item(22, 10).Value = 0.07074440112486968
$ws.Cells.Item(22, 13).Value = 5.062537867036667
$ws.Cells.Item(22, 14).Value = 1.248241426357609
$ws.Cells.Item(23, 2).Value = 0.1852027350620915
$ws.Cells.Item(23, 4).Value = 0.5312978913959796
$ws.Cells.Item(23, 5).Value = 0.12366855950917
$ws.Cells.Item(23, 6).Value = 10.28866547319222
$ws.Cells.Item(23, 7).Value = 0.002606251087721186
$ws.Cells.Item(23, 9).Value = 0.3680269456241394
$ws.Cells.Item(23, 10).Value = 0.07029105557555226
$ws.Cells.Item(23, 13).Value = 4.94430750361667
$ws.Cells.Item(23, 14).Value = 1.249960185824392
$ws.Cells.Item(24, 2).Value = 0.1605470485250891
$ws.Cells.Item(24, 4).Value = 0.4606809267102108
$ws.Cells.Item(24, 5).Value = 0.1056095816501568
$ws.Cells.Item(24, 6).Value = 9.739167168746292
$ws.Cells.Item(24, 7).Value = 0.002633723279212352
$ws.Cells.Item(24, 9).Value = 0.3796028869622337
$ws.Cells.Item(24, 10).Value = 0.06858596919020599
$ws.Cells.Item(24, 13).Value = 4.503600965836597
$ws.Cells.Item(24, 14).Value = 1.257515545999993
$ws.Cells.Item(25, 2).Value = 0.1338778525281015
$ws.Cells.Item(25, 4).Value = 0.3865000600602286
$ws.Cells.Item(25, 5).Value = 0.08634079948695472
$ws.Cells.Item(25, 6).Value = 9.186740509996582
$ws.Cells.Item(25, 7).Value = 0.002665120092791551
$ws.Cells.Item(25, 9).Value = 0.3931006435439315
$ws.Cells.Item(25, 10).Value = 0.06676978666531141
$ws.Cells.Item(25, 13).Value = 4.042081441369533
$ws.Cells.Item(25, 14).Value = 1.267923560764686
